$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Time Series" file name references (column D) to point at the
# refreshed 2024 commercial-data plots generated on 2025-05-05.
$ws.Range("D2").Value = "Commercial_LONGFINSQUID_Landings_LBS_2025-05-05.png"
$ws.Range("D3").Value = "N_Commercial_Vessels_Landing_LONGFINSQUID_2025-05-05.png"
$ws.Range("D4").Value = "TOTALANNUALREV_LONGFINSQUID_2024Dols_2025-05-05.png"

# Move the active selection to D4, matching the saved view state.
$ws.Range("D4").Select()
